$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Price" column (D) values
$ws.Range("D2").Value = "28.917.12"
$ws.Range("D3").Value = "1.823.08"
$ws.Range("D4").Value = "'0.9941"
$ws.Range("D5").Value = "'243.28"
$ws.Range("D6").Value = "'0.6288"
$ws.Range("D7").Value = "'0.9967"
$ws.Range("D8").Value = "'0.07455"
$ws.Range("D9").Value = "'0.2927"
$ws.Range("D10").Value = "'22.98"
$ws.Range("D11").Value = "'0.07672"
$ws.Range("D12").Value = "1.826.69"
$ws.Range("D14").Value = "'0.6655"
$ws.Range("D15").Value = "'82.87"
$ws.Range("D16").Value = "'0.000009698"
$ws.Range("D17").Value = "'6.015"
$ws.Range("D18").Value = "28.957.47"
$ws.Range("D20").Value = "'224.90"
$ws.Range("D21").Value = "'0.9943"
$ws.Range("D22").Value = "'7.109"
$ws.Range("D23").Value = "'0.9958"
$ws.Range("D24").Value = "'160.15"
$ws.Range("D25").Value = "'0.1406"
$ws.Range("D26").Value = "'8.483"
$ws.Range("D28").Value = "'1.494"
$ws.Range("D29").Value = "'4.107"
$ws.Range("D30").Value = "'4.044"
$ws.Range("D31").Value = "'0.05445"
$ws.Range("D32").Value = "'1.196"
$ws.Range("D34").Value = "'0.7419"
$ws.Range("D36").Value = "'2.605"
$ws.Range("D37").Value = "1.239.36"
$ws.Range("D38").Value = "'2.742"
$ws.Range("D39").Value = "'0.01775"
$ws.Range("D40").Value = "'6.650"
$ws.Range("D41").Value = "'0.8964"
$ws.Range("D42").Value = "'0.9962"
$ws.Range("D43").Value = "'101.24"
$ws.Range("D44").Value = "1.969.75"
$ws.Range("D45").Value = "'64.80"
$ws.Range("D47").Value = "'0.5060"
$ws.Range("D48").Value = "'0.4046"
$ws.Range("D49").Value = "'0.07413"
$ws.Range("D50").Value = "'8.927"
$ws.Range("D51").Value = "'1.657"

# Update "Volume(1h)" column (E) values
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("E4").Value = "  -0.50%  "
$ws.Range("E5").Value = "  +0.73%  "
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("E8").Value = "  -1.58%  "
$ws.Range("E9").Value = "  +0.27%  "
$ws.Range("E10").Value = "  +0.49%  "
$ws.Range("E11").Value = "  +0.33%  "
$ws.Range("E12").Value = "  -0.32%  "
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("E16").Value = "  +2.48%  "
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("E19").Value = "  +1.57%  "
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("E23").Value = "  -0.48%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  +2.82%  "
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("E31").Value = "  +4.47%  "
$ws.Range("E32").Value = "  -0.19%  "
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("E35").Value = "  -1.68%  "
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  -2.55%  "
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("E42").Value = "  -0.45%  "
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("E45").Value = "  +0.46%  "
$ws.Range("E46").Value = "  +1.92%  "
$ws.Range("E47").Value = "  -1.05%  "
$ws.Range("E48").Value = "  +1.68%  "
$ws.Range("E49").Value = "  +5.38%  "
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("E51").Value = "  +1.35%  "
